$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: new book entry "No Stopping You" ---
# Copy formatting from the row above first so the date cell reuses the
# existing date style instead of Excel creating a brand-new one.
$ws.Range("C15").Copy($ws.Range("C16"))
$ws.Range("D15").Copy($ws.Range("D16"))

$ws.Range("A16").Value = "No Stopping You"
$ws.Range("B16").Value = "Roger Flax"
$ws.Range("C16").Value = "2/13/2021"
$ws.Range("D16").Value = "2/21/2021"
$ws.Range("E16").Value = "self improvement;business;public speaking;success"
$ws.Range("F16").Value = "Audio"
$ws.Range("G16").Value = "11 Hours 4 Mins"
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = $true

# --- Row 17: new (in-progress) book entry "Harry Potter and the Deathly Hallows" ---
$ws.Range("C15").Copy($ws.Range("C17"))

$ws.Range("A17").Value = "Harry Potter and the Deathly Hallows"
$ws.Range("B17").Value = "J.K. Rowling"
$ws.Range("C17").Value = "2/19/2021"

$ws.Range("E17").Select()
